$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# 1) Text change: "Ready for handoff" -> "In Translation" wherever it
#    appears (Status columns on the Overview/zh-cn/de-de sheets).
# -----------------------------------------------------------------
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        $v = $cell.Value2()
        if ($v -is [string] -and "Ready for handoff" -eq $v) {
            $cell.Value2 = "In Translation"
        }
    }
}

# -----------------------------------------------------------------
# 2) Column widths shrink to follow the shorter status text
#    (re-fit of the "Status"/"zh-cn"/"de-de" columns).
# -----------------------------------------------------------------
$newWidth = 13.4101845877511

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = $newWidth   # column E (zh-cn)
$wsOverview.Columns.Item(6).ColumnWidth = $newWidth   # column F (de-de)

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).ColumnWidth = $newWidth        # column C (Status)

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).ColumnWidth = $newWidth        # column C (Status)
